$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by Excel (e.g. "676.43")
# are forced to Text format first so they are stored as strings, matching the
# original workbook where every Price/Volume cell is an inline/shared string.
$textFormatCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped from the latest cryptos snapshot.
$ws.Range('D2').Value = '68.942.15'
$ws.Range('E2').Value = '  -2.90%  '
$ws.Range('D3').Value = '3.667.03'
$ws.Range('E3').Value = '  -4.22%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '676.43'
$ws.Range('E5').Value = '  -4.28%  '
$ws.Range('D6').Value = '160.54'
$ws.Range('D7').Value = '3.665.40'
$ws.Range('E7').Value = '  -4.24%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').Value = '0.488'
$ws.Range('E9').Value = '  -6.83%  '
$ws.Range('D10').Value = '0.146'
$ws.Range('E10').Value = '  -9.69%  '
$ws.Range('D11').Value = '7.18'
$ws.Range('E11').Value = '  -2.90%  '
$ws.Range('D12').Value = '0.444'
$ws.Range('E12').Value = '  -3.13%  '
$ws.Range('D13').Value = '0.0000232'
$ws.Range('E13').Value = '  -8.95%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '4.299.00'
$ws.Range('E14').Value = '  -4.01%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '32.94'
$ws.Range('E15').Value = '  -9.77%  '
$ws.Range('D16').Value = '3.671.47'
$ws.Range('E16').Value = '  -4.56%  '
$ws.Range('D17').Value = '68.959.85'
$ws.Range('E17').Value = '  -2.98%  '
$ws.Range('E18').Value = '  -1.82%  '
$ws.Range('D19').Value = '16.10'
$ws.Range('E19').Value = '  -7.27%  '
$ws.Range('D20').Value = '6.49'
$ws.Range('E20').Value = '  -9.92%  '
$ws.Range('D21').Value = '477.05'
$ws.Range('E21').Value = '  -3.48%  '
$ws.Range('D22').Value = '9.71'
$ws.Range('E22').Value = '  -8.59%  '
$ws.Range('D23').Value = '0.657'
$ws.Range('E23').Value = '  -10.45%  '
$ws.Range('D24').Value = '78.74'
$ws.Range('E24').Value = '  -8.08%  '
$ws.Range('D25').Value = '3.821.96'
$ws.Range('E25').Value = '  -4.03%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '11.43'
$ws.Range('E27').Value = '  -5.61%  '
$ws.Range('D28').Value = '0.0000125'
$ws.Range('E28').Value = '  -13.42%  '
$ws.Range('D29').Value = '9.23'
$ws.Range('E29').Value = '  -12.89%  '
$ws.Range('D30').Value = '1.78'
$ws.Range('E30').Value = '  -14.35%  '
$ws.Range('D31').Value = '2.70'
$ws.Range('E31').Value = '  -12.71%  '
$ws.Range('D32').Value = '2.08'
$ws.Range('E32').Value = '  -7.16%  '
$ws.Range('D33').Value = '6.63'
$ws.Range('E33').Value = '  -10.44%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.166'
$ws.Range('E35').Value = '  -5.62%  '
$ws.Range('D36').Value = '26.54'
$ws.Range('E36').Value = '  -9.46%  '
$ws.Range('D37').Value = '3.642.93'
$ws.Range('E37').Value = '  -4.09%  '
$ws.Range('D38').Value = '8.40'
$ws.Range('E38').Value = '  -8.17%  '
$ws.Range('D39').Value = '5.96'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').Value = '0.0920'
$ws.Range('E40').Value = '  -9.79%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '2.15'
$ws.Range('E43').Value = '  -7.31%  '
$ws.Range('E44').Value = '  -9.73%  '
$ws.Range('D45').Value = '160.31'
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('D46').Value = '48.34'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('D47').Value = '2.80'
$ws.Range('E47').Value = '  -15.49%  '
$ws.Range('D48').Value = '1.29'
$ws.Range('E48').Value = '  -5.17%  '
$ws.Range('D49').Value = '381.54'
$ws.Range('E49').Value = '  -11.22%  '
$ws.Range('D50').Value = '0.000270'
$ws.Range('E50').Value = '  -12.96%  '
$ws.Range('D51').Value = '7.94'
$ws.Range('E51').Value = '  -9.37%  '
